# Cover letter edits per commit "cover letter .docx version."

$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "With the TensorFlow Lite models integrated, ..." -> bold "TensorFlow Lite models"
# Locate the phrase to bold and apply bold to just that sub-range; Word
# automatically splits the surrounding run(s) to preserve the rest of the
# paragraph's (non-bold) formatting.
$r1b = $d.Content
$foundBold = $r1b.Find.Execute("TensorFlow Lite models")
if ($foundBold) {
    $r1b.Bold = 1
}

# --- Change 2 -------------------------------------------------------------
# "Detect and immediately send payload notification to the car owner
#  confirming the car is being driven by an unrecognized driver."
# ->
# "Detect and immediately sends payload notification to car owner
#  confirming the car is being driven by an unrecognized usual driver."
$d.Content.Find.Execute(
    "Detect and immediately send payload notification to the car owner confirming the car is being driven by an unrecognized driver. If owner confirms ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Detect and immediately sends payload notification to car owner confirming the car is being driven by an unrecognized usual driver. If owner confirms ",
    2
) | Out-Null

# Bold run "“Block!”" -> "“REPORT”"
$d.Content.Find.Execute(
    [char]0x201C + "Block!" + [char]0x201D,
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x201C + "REPORT" + [char]0x201D,
    2
) | Out-Null

# --- Change 3 -------------------------------------------------------------
# "... with local models / or downloaded from Firebase." ->
# "... with local models / or remote downloaded from Firebase."
$d.Content.Find.Execute(
    "local models / or downloaded from Firebase.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "local models / or remote downloaded from Firebase.",
    2
) | Out-Null
